$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.311.46"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").Value = "'2.580.06"
$ws.Range("E3").Value = "  -0.90%  "

$ws.Range("D5").Value = "'573.35"
$ws.Range("E5").Value = "  +2.99%  "

$ws.Range("D6").Value = "'143.27"
$ws.Range("E6").Value = "  -0.31%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").Value = "'2.587.32"
$ws.Range("E9").Value = "  -1.11%  "

$ws.Range("E10").Value = "  -1.60%  "

$ws.Range("E11").Value = "  +2.96%  "

$ws.Range("E12").Value = "  +11.19%  "

$ws.Range("E13").Value = "  +3.13%  "

$ws.Range("D14").Value = "'3.034.86"
$ws.Range("E14").Value = "  -1.04%  "

$ws.Range("D15").Value = "'59.312.21"
$ws.Range("E15").Value = "  +0.78%  "

$ws.Range("D16").Value = "'22.47"
$ws.Range("E16").Value = "  +7.14%  "

$ws.Range("D17").Value = "'0.0000137"
$ws.Range("E17").Value = "  +3.87%  "

$ws.Range("D18").Value = "'2.584.40"
$ws.Range("E18").Value = "  -1.14%  "

$ws.Range("E19").Value = "  +1.55%  "

$ws.Range("D20").Value = "'339.77"
$ws.Range("E20").Value = "  +0.64%  "

$ws.Range("E21").Value = "  +1.73%  "

$ws.Range("D22").Value = "'6.28"
$ws.Range("E22").Value = "  +1.73%  "

$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").Value = "'64.66"
$ws.Range("E24").Value = "  -2.45%  "

$ws.Range("E25").Value = "  +7.48%  "

$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("E27").Value = "  -0.31%  "

$ws.Range("D28").Value = "'7.29"
$ws.Range("E28").Value = "  +1.95%  "

$ws.Range("D29").Value = "'0.0₃0783"
$ws.Range("E29").Value = "  +3.02%  "

$ws.Range("D32").Value = "'6.07"
$ws.Range("E32").Value = "  +1.21%  "

$ws.Range("D33").Value = "'158.07"
$ws.Range("E33").Value = "  +2.57%  "

$ws.Range("D34").Value = "'19.05"
$ws.Range("E34").Value = "  +0.62%  "

$ws.Range("E35").Value = "  +1.97%  "

$ws.Range("E36").Value = "  +2.42%  "

$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "'0.879"
$ws.Range("E37").Value = "  -3.07%  "

$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").Value = "'0.872"
$ws.Range("E38").Value = "  -3.86%  "

$ws.Range("D39").Value = "'37.27"
$ws.Range("E39").Value = "  +0.40%  "

$ws.Range("E40").Value = "  +1.75%  "

$ws.Range("D41").Value = "'296.88"
$ws.Range("E41").Value = "  +4.97%  "

$ws.Range("D42").Value = "'3.68"
$ws.Range("E42").Value = "  +2.56%  "

$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  +0.21%  "

$ws.Range("D44").Value = "'0.0978"
$ws.Range("E44").Value = "  +2.52%  "

$ws.Range("D45").Value = "'129.75"
$ws.Range("E45").Value = "  +9.93%  "

$ws.Range("D46").Value = "'0.595"
$ws.Range("E46").Value = "  -1.08%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'19.25"
$ws.Range("E47").Value = "  +2.56%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0537"
$ws.Range("E48").Value = "  -0.27%  "

$ws.Range("D49").Value = "'10.66"
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("E50").Value = "  +2.86%  "

$ws.Range("D51").Value = "'1.956.20"
$ws.Range("E51").Value = "  +0.46%  "
